# Apply the "合肥-漫展信息" update:
#  - 展览   (Exhibitions) sheet: bump several "want to go" counts, and insert a
#            brand-new exhibition row ("合肥·第十五届次元之门动漫游戏博览会")
#            before the existing "首届AT次元时代" row, shifting later rows down.
#  - 演出   (Performances) sheet: bump two "want to go" counts.
#  - 全部类型 (All types) sheet: same bumps + same new-row insertion as 展览,
#            since it aggregates every row from the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

# Simple "want to go" count updates for existing rows.
$ws.Range("F2").Value = 184
$ws.Range("F3").Value = 501
$ws.Range("F4").Value = 26
$ws.Range("F5").Value = 22
$ws.Range("F8").Value = 22

# Insert a new row at position 9; this pushes the old rows 9 and 10 down to
# rows 10 and 11 (together with their formatting and content).
$ws.Rows.Item(9).Insert()

# Copy the formatting of the (now shifted) row 10 into the freshly inserted
# row 9 so the new row matches the sheet's existing look (bold/bordered id
# column, etc.).
$ws.Range("A10:I10").Copy($ws.Range("A9:I9"))

# Fill in the new exhibition's data.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "2024-10-01"
$ws.Range("C9").Value = "合肥·第十五届次元之门动漫游戏博览会"
$ws.Range("D9").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Range("E9").Value = "2024.10.01 09:30-10.02 17:30"
$ws.Range("F9").Value = 60
$ws.Range("G9").Value = 70
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=91133"
$ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg"

# Update the "want to go" counts of the two rows that got shifted down.
$ws.Range("F10").Value = 1802
$ws.Range("F11").Value = 5

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")

$ws.Range("F2").Value = 92
$ws.Range("F3").Value = 41

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - mirrors the 展览 sheet changes, offset by the
# extra rows already present (演出 rows) before the inserted exhibition.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F2").Value = 184
$ws.Range("F3").Value = 92
$ws.Range("F4").Value = 501
$ws.Range("F5").Value = 26
$ws.Range("F6").Value = 22
$ws.Range("F9").Value = 22

# Insert the same new exhibition row, this time at position 10.
$ws.Rows.Item(10).Insert()
$ws.Range("A11:I11").Copy($ws.Range("A10:I10"))

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "2024-10-01"
$ws.Range("C10").Value = "合肥·第十五届次元之门动漫游戏博览会"
$ws.Range("D10").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Range("E10").Value = "2024.10.01 09:30-10.02 17:30"
$ws.Range("F10").Value = 60
$ws.Range("G10").Value = 70
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=91133"
$ws.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg"

# Update "want to go" counts for rows shifted down by the insertion.
$ws.Range("F11").Value = 1802
$ws.Range("F12").Value = 5
$ws.Range("F13").Value = 41
